$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '68.954.00'
$ws.Range('E2').Value = '  -2.66%  '
$ws.Range('D3').Value = '3.518.98'
$ws.Range('E3').Value = '  -3.28%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '587.45'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.55%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '169.90'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -3.61%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.613'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.39%  '
$ws.Range('D8').Value = '3.510.20'
$ws.Range('E8').Value = '  -3.36%  '
$ws.Range('E9').Value = '  +0.02%  '
$ws.Range('E10').Value = '  -4.50%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.77'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -1.12%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.577'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -5.12%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '47.28'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -2.69%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000275'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -3.41%  '
$ws.Range('D15').Value = '4.082.67'
$ws.Range('E15').Value = '  -3.35%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '8.42'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -6.01%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '612.12'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -9.47%  '
$ws.Range('D18').Value = '68.969.08'
$ws.Range('E18').Value = '  -2.64%  '
$ws.Range('D19').Value = '3.510.25'
$ws.Range('E19').Value = '  -3.33%  '
$ws.Range('E20').Value = '  -1.10%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.38'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -2.46%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.09'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -3.71%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.884'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -6.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '15.72'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -8.51%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '96.40'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -3.68%  '
$ws.Range('E26').Value = '  -2.38%  '
$ws.Range('E27').Value = '  +0.01%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.61'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -6.85%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.19'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -6.68%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '32.53'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -6.25%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.50'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -7.14%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.12'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -5.63%  '
$ws.Range('E33').Value = '  -5.15%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.89'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -9.18%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '613.61'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +6.45%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '10.71'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -3.53%  '
$ws.Range('B37').Value = 'dogwifhat'
$ws.Range('C37').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.45'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -13.86%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.102'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -4.97%  '
$ws.Range('B39').Value = 'OKB'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '57.09'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -2.55%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.999'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +0.04%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0441'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -3.05%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.136'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -3.36%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '3.381.11'
$ws.Range('E43').Value = '  -5.25%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.325'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -6.22%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '32.67'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -5.31%  '
$ws.Range('E46').Value = '  -5.38%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.51'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -6.73%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.74'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -5.49%  '
$ws.Range('E49').Value = '  -3.42%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '133.53'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -2.79%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '5.58'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +11.59%  '
